# Fruta / hortaliza, semanal
# Insert two new data rows (new rows 51 and 52) above the current row 51,
# shifting the existing rows 51-77 down to 53-79 (Excel COM Insert with
# shift-down, matching xlShiftDown = -4121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51:52").Insert(-4121)

# Common (constant across the whole "Membrillo" subset) column values.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100104
$producto   = "Frutos de pepita"
$categoriaId = 100104003
$categoria   = "Membrillo"
$variedad    = "Champion"

# New row 51
$ws.Range("A51").Value = $mercadoId
$ws.Range("B51").Value = $mercado
$ws.Range("C51").Value = $region
$ws.Range("D51").Value = 45068
$ws.Range("E51").Value = $codreg
$ws.Range("F51").Value = $tipo
$ws.Range("G51").Value = $productoId
$ws.Range("H51").Value = $producto
$ws.Range("I51").Value = $categoriaId
$ws.Range("J51").Value = $categoria
$ws.Range("K51").Value = $variedad
$ws.Range("L51").Value = "Especial"
$ws.Range("M51").Value = 220
$ws.Range("N51").Value = 10500
$ws.Range("O51").Value = 10500
$ws.Range("P51").Value = 10500
$ws.Range("Q51").Value = '$/caja 15 kilos empedrada'
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 700
$ws.Range("T51").Value = 15

# New row 52
$ws.Range("A52").Value = $mercadoId
$ws.Range("B52").Value = $mercado
$ws.Range("C52").Value = $region
$ws.Range("D52").Value = 45068
$ws.Range("E52").Value = $codreg
$ws.Range("F52").Value = $tipo
$ws.Range("G52").Value = $productoId
$ws.Range("H52").Value = $producto
$ws.Range("I52").Value = $categoriaId
$ws.Range("J52").Value = $categoria
$ws.Range("K52").Value = $variedad
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 200
$ws.Range("N52").Value = 7500
$ws.Range("O52").Value = 7500
$ws.Range("P52").Value = 7500
$ws.Range("Q52").Value = '$/caja 15 kilos empedrada'
$ws.Range("R52").Value = "Región de O'Higgins"
$ws.Range("S52").Value = 500
$ws.Range("T52").Value = 15
